$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.458825068571144
$ws.Range("C2").Value = 0.04077326317384689
$ws.Range("E2").Value = 0.04600774054123669
$ws.Range("F2").Value = 6.354418123488841
$ws.Range("G2").Value = 0.002694860732103011
$ws.Range("J2").Value = 0.6049521132995181
$ws.Range("K2").Value = 1.062898978091397
$ws.Range("L2").Value = 0.1207030533839752
$ws.Range("M2").Value = 0.2533381479270957
$ws.Range("B3").Value = 1.454157706160146
$ws.Range("C3").Value = 0.03551341101903915
$ws.Range("E3").Value = 0.04585883875205532
$ws.Range("F3").Value = 6.155580197608202
$ws.Range("G3").Value = 0.002699896048564795
$ws.Range("J3").Value = 0.579872450894058
$ws.Range("K3").Value = 1.04783688151565
$ws.Range("L3").Value = 0.1246140913417033
$ws.Range("M3").Value = 0.2571453272431832
$ws.Range("B4").Value = 1.452452678315041
$ws.Range("C4").Value = 0.03228020722961844
$ws.Range("E4").Value = 0.04580120698253509
$ws.Range("F4").Value = 6.03387621080708
$ws.Range("G4").Value = 0.00270314829939684
$ws.Range("J4").Value = 0.5643051192804336
$ws.Range("K4").Value = 1.039458673152325
$ws.Range("L4").Value = 0.1271880884942007
$ws.Range("M4").Value = 0.2597445266767444
$ws.Range("B5").Value = 1.452049623596679
$ws.Range("C5").Value = 0.03096163090714299
$ws.Range("E5").Value = 0.04578619029921072
$ws.Range("F5").Value = 5.984373130962865
$ws.Range("G5").Value = 0.002704514130386054
$ws.Range("J5").Value = 0.5579181552496806
$ws.Range("K5").Value = 1.036262914896582
$ws.Range("L5").Value = 0.1282803239314649
$ws.Range("M5").Value = 0.260869483025056
$ws.Range("B6").Value = 1.452000314847169
$ws.Range("C6").Value = 0.03074261766168718
$ws.Range("E6").Value = 0.04578420743186129
$ws.Range("F6").Value = 5.976158636998434
$ws.Range("G6").Value = 0.002704743376639617
$ws.Range("J6").Value = 0.5568549745120634
$ws.Range("K6").Value = 1.035745444799915
$ws.Range("L6").Value = 0.1284643019821559
$ws.Range("M6").Value = 0.2610602526414034
$ws.Range("B7").Value = 1.45244606142623
$ws.Range("C7").Value = 0.03226242871359375
$ws.Range("E7").Value = 0.04580097021128893
$ws.Range("F7").Value = 6.033208226744563
$ws.Range("G7").Value = 0.002703166555243872
$ws.Range("J7").Value = 0.5642191582385152
$ws.Range("K7").Value = 1.039414690065882
$ws.Range("L7").Value = 0.127202643523038
$ws.Range("M7").Value = 0.2597594319860281
$ws.Range("B8").Value = 1.45697476281066
$ws.Range("C8").Value = 0.03896035481263027
$ws.Range("E8").Value = 0.04594936734192068
$ws.Range("F8").Value = 6.285776712741125
$ws.Range("G8").Value = 0.002696563673350771
$ws.Range("J8").Value = 0.5963390770466361
$ws.Range("K8").Value = 1.057524789508449
$ws.Range("L8").Value = 0.1220156999448108
$ws.Range("M8").Value = 0.2545965639211545
$ws.Range("B9").Value = 1.475075218608623
$ws.Range("C9").Value = 0.05207153765013572
$ws.Range("E9").Value = 0.04650992771597018
$ws.Range("F9").Value = 6.784314538435865
$ws.Range("G9").Value = 0.002684882834731794
$ws.Range("J9").Value = 0.6580301107113797
$ws.Range("K9").Value = 1.099961077842636
$ws.Range("L9").Value = 0.1132177493942592
$ws.Range("M9").Value = 0.246548609612038
$ws.Range("B10").Value = 1.494013532783669
$ws.Range("C10").Value = 0.06169841116089003
$ws.Range("E10").Value = 0.04708818649750768
$ws.Range("F10").Value = 7.152900717705194
$ws.Range("G10").Value = 0.002677064481772184
$ws.Range("J10").Value = 0.7026223693768543
$ws.Range("K10").Value = 1.135392672235383
$ws.Range("L10").Value = 0.1075968941976804
$ws.Range("M10").Value = 0.2419033400470525
$ws.Range("B11").Value = 1.503858533986005
$ws.Range("C11").Value = 0.06607851215122196
$ws.Range("E11").Value = 0.04738785016007441
$ws.Range("F11").Value = 7.321155453849542
$ws.Range("G11").Value = 0.002673671569533908
$ws.Range("J11").Value = 0.7227622763728334
$ws.Range("K11").Value = 1.152442933847738
$ws.Range("L11").Value = 0.1052241060798096
$ws.Range("M11").Value = 0.2400657818398635
$ws.Range("B12").Value = 1.507763755512372
$ws.Range("C12").Value = 0.06773738944980323
$ws.Range("E12").Value = 0.04750662539372641
$ws.Range("F12").Value = 7.384958266694923
$ws.Range("G12").Value = 0.002672410153066098
$ws.Range("J12").Value = 0.7303687952854148
$ws.Range("K12").Value = 1.159034028458564
$ws.Range("L12").Value = 0.1043521857423535
$ws.Range("M12").Value = 0.2394096177018454
$ws.Range("B13").Value = 1.506914813492472
$ws.Range("C13").Value = 0.06738010855437437
$ws.Range("E13").Value = 0.04748080878887251
$ws.Range("F13").Value = 7.371213184655232
$ws.Range("G13").Value = 0.002672680782596009
$ws.Range("J13").Value = 0.7287314720881
$ws.Range("K13").Value = 1.157608527299715
$ws.Range("L13").Value = 0.1045387842767056
$ws.Range("M13").Value = 0.2395491688340563
$ws.Range("B14").Value = 1.504176267233674
$ws.Range("C14").Value = 0.06621498375113788
$ws.Range("E14").Value = 0.04739751546881266
$ws.Range("F14").Value = 7.326402751739352
$ws.Range("G14").Value = 0.002673567323561971
$ws.Range("J14").Value = 0.7233884656099292
$ws.Range("K14").Value = 1.152982488203321
$ws.Range("L14").Value = 0.1051518388102153
$ws.Range("M14").Value = 0.2400110033000011
$ws.Range("B15").Value = 1.502521905149649
$ws.Range("C15").Value = 0.06550134430453625
$ws.Range("E15").Value = 0.04734718703672414
$ws.Range("F15").Value = 7.298966715829408
$ws.Range("G15").Value = 0.002674113399746769
$ws.Range("J15").Value = 0.720113137778668
$ws.Range("K15").Value = 1.150166439876244
$ws.Range("L15").Value = 0.1055308205103955
$ws.Range("M15").Value = 0.240299059154804
$ws.Range("B16").Value = 1.493394910402969
$ws.Range("C16").Value = 0.0614121868987354
$ws.Range("E16").Value = 0.04706934241734473
$ws.Range("F16").Value = 7.141916995088508
$ws.Range("G16").Value = 0.002677289498960257
$ws.Range("J16").Value = 0.7013033290062936
$ws.Range("K16").Value = 1.134297196971346
$ws.Range("L16").Value = 0.1077556776064625
$ws.Range("M16").Value = 0.2420289786183183
$ws.Range("B17").Value = 1.488110986902427
$ws.Range("C17").Value = 0.05890390435075687
$ws.Range("E17").Value = 0.04690829508954764
$ws.Range("F17").Value = 7.045724628025027
$ws.Range("G17").Value = 0.002679279764023673
$ws.Range("J17").Value = 0.6897275155741625
$ws.Range("K17").Value = 1.124801034561898
$ws.Range("L17").Value = 0.1091678097168121
$ws.Range("M17").Value = 0.2431608490453101
$ws.Range("B18").Value = 1.485187557013433
$ws.Range("C18").Value = 0.05746127378823473
$ws.Range("E18").Value = 0.04681910882552742
$ws.Range("F18").Value = 6.99045192412629
$ws.Range("G18").Value = 0.002680439927710885
$ws.Range("J18").Value = 0.6830556817217541
$ws.Range("K18").Value = 1.119426811734911
$ws.Range("L18").Value = 0.1099973581390934
$ws.Range("M18").Value = 0.2438378084919748
$ws.Range("B19").Value = 1.484217604621421
$ws.Range("C19").Value = 0.0569728319380971
$ws.Range("E19").Value = 0.04678950233648749
$ws.Range("F19").Value = 6.971746776414562
$ws.Range("G19").Value = 0.002680835390907963
$ws.Range("J19").Value = 0.6807943284356384
$ws.Range("K19").Value = 1.117622243799417
$ws.Range("L19").Value = 0.1102812013506345
$ws.Range("M19").Value = 0.244071468689981
$ws.Range("B20").Value = 1.488661488661592
$ws.Range("C20").Value = 0.05917090752338083
$ws.Range("E20").Value = 0.04692508217372904
$ws.Range("F20").Value = 7.055958785912821
$ws.Range("G20").Value = 0.002679066302489643
$ws.Range("J20").Value = 0.6909611955699972
$ws.Range("K20").Value = 1.125802835354364
$ws.Range("L20").Value = 0.1090156916674481
$ws.Range("M20").Value = 0.2430376747390071
$ws.Range("B21").Value = 1.504975834992337
$ws.Range("C21").Value = 0.06655720203535509
$ws.Range("E21").Value = 0.0474218366426129
$ws.Range("F21").Value = 7.339562226420071
$ws.Range("G21").Value = 0.0026733062910056
$ws.Range("J21").Value = 0.7249583726703008
$ws.Range("K21").Value = 1.154337613043793
$ws.Range("L21").Value = 0.1049710468672274
$ws.Range("M21").Value = 0.239874274034829
$ws.Range("B22").Value = 1.516670813653889
$ws.Range("C22").Value = 0.07138595637675849
$ws.Range("E22").Value = 0.04777739797071234
$ws.Range("F22").Value = 7.525430871092397
$ws.Range("G22").Value = 0.002669678154485045
$ws.Range("J22").Value = 0.7470611708787374
$ws.Range("K22").Value = 1.173771118512917
$ws.Range("L22").Value = 0.10248275769084
$ws.Range("M22").Value = 0.2380381148347475
$ws.Range("B23").Value = 1.510334402570749
$ws.Range("C23").Value = 0.06880859438925313
$ws.Range("E23").Value = 0.0475847890101484
$ws.Range("F23").Value = 7.426180411999951
$ws.Range("G23").Value = 0.002671602125841447
$ws.Range("J23").Value = 0.7352748452721869
$ws.Range("K23").Value = 1.163327159367498
$ws.Range("L23").Value = 0.1037965707499069
$ws.Range("M23").Value = 0.2389969261421214
$ws.Range("B24").Value = 1.488412250590756
$ws.Range("C24").Value = 0.05905019722386839
$ws.Range("E24").Value = 0.04691748213901903
$ws.Range("F24").Value = 7.05133183135473
$ws.Range("G24").Value = 0.002679162758837517
$ws.Range("J24").Value = 0.6904035010804535
$ws.Range("K24").Value = 1.125349655715439
$ws.Range("L24").Value = 0.1090844091286876
$ws.Range("M24").Value = 0.2430932801350991
$ws.Range("B25").Value = 1.469189595313537
$ws.Range("C25").Value = 0.04852626294901086
$ws.Range("E25").Value = 0.04632920912608718
$ws.Range("F25").Value = 6.649065039614925
$ws.Range("G25").Value = 0.002687908068522038
$ws.Range("J25").Value = 0.6414738285291008
$ws.Range("K25").Value = 1.087736430820343
$ws.Range("L25").Value = 0.1154501998614919
$ws.Range("M25").Value = 0.248503332835984
